# This script applies updated calibration sampling ranges (max_35 / min_35)
# to the IPPU input-variable sheets, and normalizes the gas-recovery-fraction
# (gasrf) trajectories for CO2 capture to a flat 0.9 factor.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("strategy_id-0")
foreach ($r in @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23)) {
    $ws.Cells.Item($r, 8).Value = 1.5   # max_35
    $ws.Cells.Item($r, 9).Value = 0.5   # min_35
}
foreach ($r in @(24, 26, 27, 28, 30, 31, 34, 35, 36, 37, 38, 39, 40, 42, 44, 45, 46, 49, 50, 51, 54, 55, 57, 59, 61, 65, 66, 69, 70, 71, 72, 73, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 104, 105, 106, 107, 108, 109, 110, 111, 113, 114, 127, 128, 129, 130, 131, 150, 151, 152, 153, 154, 155, 156, 157, 158, 159, 160, 161, 162, 163, 164, 169, 170, 171, 172, 173, 174, 175, 176, 177, 178, 179, 180, 181, 182, 183, 184, 185, 186)) {
    $ws.Cells.Item($r, 8).Value = 1   # max_35
    $ws.Cells.Item($r, 9).Value = 1   # min_35
}
foreach ($r in @(122, 123, 124, 125, 126)) {
    $ws.Range("J" + $r + ":AS" + $r).Value = 0.9
}

$ws = $wb.Worksheets.Item("strategy_id-5001")
foreach ($r in @(2)) {
    $ws.Cells.Item($r, 8).Value = 1   # max_35
    $ws.Cells.Item($r, 9).Value = 1   # min_35
}

$ws = $wb.Worksheets.Item("strategy_id-5006")
foreach ($r in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)) {
    $ws.Cells.Item($r, 8).Value = 1.5   # max_35
    $ws.Cells.Item($r, 9).Value = 0.5   # min_35
}
foreach ($r in @(22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34)) {
    $ws.Cells.Item($r, 8).Value = 1   # max_35
    $ws.Cells.Item($r, 9).Value = 1   # min_35
}

$ws = $wb.Worksheets.Item("strategy_id-5008")
foreach ($r in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)) {
    $ws.Cells.Item($r, 8).Value = 1.5   # max_35
    $ws.Cells.Item($r, 9).Value = 0.5   # min_35
}
foreach ($r in @(22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35)) {
    $ws.Cells.Item($r, 8).Value = 1   # max_35
    $ws.Cells.Item($r, 9).Value = 1   # min_35
}
